# Updates odds data (columns G:AJ) on Sheet1 for Jogos_da_Semana_FlashScore_2025-06-04
# to match the latest FlashScore scrape, per commit "Atualizando o arquivo XLSX".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("G6").Value = 1.62
$ws.Range("H6").Value = 3.6
$ws.Range("I6").Value = 5.2
$ws.Range("L6").Value = 1.33
$ws.Range("M6").Value = 2.8
$ws.Range("N6").Value = 1.98
$ws.Range("O6").Value = 1.65
$ws.Range("P6").Value = 1.44
$ws.Range("Q6").Value = 2.42
$ws.Range("R6").Value = 1.98
$ws.Range("S6").Value = 1.65
$ws.Range("T6").Value = 5.8
$ws.Range("U6").Value = 6.8
$ws.Range("V6").Value = 8.25
$ws.Range("W6").Value = 11.5
$ws.Range("X6").Value = 14
$ws.Range("Y6").Value = 32
$ws.Range("Z6").Value = 8.75
$ws.Range("AA6").Value = 7.1
$ws.Range("AB6").Value = 18.5
$ws.Range("AC6").Value = 100
$ws.Range("AD6").Value = 1000
$ws.Range("AE6").Value = 12.5
$ws.Range("AF6").Value = 30
$ws.Range("AG6").Value = 17
$ws.Range("AH6").Value = 100
$ws.Range("AI6").Value = 60
$ws.Range("AJ6").Value = 65
# Row 7
$ws.Range("G7").Value = 1.37
$ws.Range("H7").Value = 4.05
$ws.Range("I7").Value = 8.75
$ws.Range("L7").Value = 1.28
$ws.Range("M7").Value = 3.05
$ws.Range("N7").Value = 1.82
$ws.Range("O7").Value = 1.8
$ws.Range("P7").Value = 1.39
$ws.Range("Q7").Value = 2.55
$ws.Range("R7").Value = 2.05
$ws.Range("S7").Value = 1.6
$ws.Range("T7").Value = 5.9
$ws.Range("U7").Value = 6
$ws.Range("V7").Value = 8.25
$ws.Range("W7").Value = 8.75
$ws.Range("X7").Value = 11.75
$ws.Range("Y7").Value = 32
$ws.Range("Z7").Value = 9.25
$ws.Range("AA7").Value = 8.25
$ws.Range("AB7").Value = 23
$ws.Range("AC7").Value = 120
$ws.Range("AE7").Value = 19
$ws.Range("AF7").Value = 65
$ws.Range("AG7").Value = 29
$ws.Range("AH7").Value = 300
$ws.Range("AI7").Value = 120
$ws.Range("AJ7").Value = 120
# Row 8
$ws.Range("G8").Value = 2.67
$ws.Range("H8").Value = 3.05
$ws.Range("I8").Value = 2.6
$ws.Range("L8").Value = 1.37
$ws.Range("M8").Value = 2.65
$ws.Range("N8").Value = 2.07
$ws.Range("O8").Value = 1.6
$ws.Range("P8").Value = 1.42
$ws.Range("Q8").Value = 2.45
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8
$ws.Range("T8").Value = 7.7
$ws.Range("U8").Value = 13
$ws.Range("V8").Value = 10
$ws.Range("W8").Value = 32
$ws.Range("X8").Value = 24
$ws.Range("Y8").Value = 35
$ws.Range("Z8").Value = 7.9
$ws.Range("AA8").Value = 5.9
$ws.Range("AB8").Value = 14.5
$ws.Range("AC8").Value = 75
$ws.Range("AD8").Value = 700
$ws.Range("AE8").Value = 7.4
$ws.Range("AF8").Value = 12.5
$ws.Range("AG8").Value = 9.75
$ws.Range("AH8").Value = 29
$ws.Range("AI8").Value = 24
$ws.Range("AJ8").Value = 35
# Row 9
$ws.Range("N9").Value = 1.34
$ws.Range("S9").Value = 2.67
$ws.Range("U9").Value = 14.5
$ws.Range("AC9").Value = 29
$ws.Range("AE9").Value = 18.5
$ws.Range("AF9").Value = 23
# Row 10
$ws.Range("G10").Value = 2.25
$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 2.62
$ws.Range("T10").Value = 11.5
$ws.Range("U10").Value = 13.5
$ws.Range("V10").Value = 9.25
$ws.Range("W10").Value = 23
$ws.Range("X10").Value = 16
$ws.Range("Y10").Value = 21
$ws.Range("AD10").Value = 200
$ws.Range("AE10").Value = 12.5
$ws.Range("AF10").Value = 16
$ws.Range("AG10").Value = 10.25
$ws.Range("AH10").Value = 30
$ws.Range("AI10").Value = 19
$ws.Range("AJ10").Value = 23
# Row 13
$ws.Range("G13").Value = 3.45
$ws.Range("H13").Value = 3.15
$ws.Range("I13").Value = 2.07
$ws.Range("L13").Value = 1.44
$ws.Range("M13").Value = 2.4
$ws.Range("N13").Value = 2.27
$ws.Range("O13").Value = 1.5
$ws.Range("P13").Value = 1.47
$ws.Range("Q13").Value = 2.32
$ws.Range("T13").Value = 7.9
$ws.Range("U13").Value = 16
$ws.Range("V13").Value = 13
$ws.Range("W13").Value = 50
$ws.Range("X13").Value = 37
$ws.Range("Y13").Value = 55
$ws.Range("Z13").Value = 7.1
$ws.Range("AA13").Value = 6.3
$ws.Range("AB13").Value = 19
$ws.Range("AC13").Value = 120
$ws.Range("AE13").Value = 5.8
$ws.Range("AF13").Value = 8.5
$ws.Range("AG13").Value = 9.25
$ws.Range("AH13").Value = 18.5
$ws.Range("AI13").Value = 20
# Row 14
$ws.Range("G14").Value = 3.55
$ws.Range("I14").Value = 2.1
$ws.Range("L14").Value = 1.47
$ws.Range("M14").Value = 2.32
$ws.Range("N14").Value = 2.37
$ws.Range("O14").Value = 1.45
$ws.Range("Q14").Value = 2.2
$ws.Range("R14").Value = 2.07
$ws.Range("S14").Value = 1.6
$ws.Range("T14").Value = 7.8
$ws.Range("U14").Value = 17
$ws.Range("V14").Value = 13.5
$ws.Range("W14").Value = 55
$ws.Range("X14").Value = 45
$ws.Range("Z14").Value = 6.6
$ws.Range("AB14").Value = 19.5
$ws.Range("AC14").Value = 120
$ws.Range("AF14").Value = 8.75
$ws.Range("AG14").Value = 9.25
$ws.Range("AH14").Value = 19
$ws.Range("AI14").Value = 21
$ws.Range("AJ14").Value = 40
# Row 15
$ws.Range("G15").Value = 2.07
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = 3.65
$ws.Range("L15").Value = 1.47
$ws.Range("M15").Value = 2.35
$ws.Range("N15").Value = 2.32
$ws.Range("O15").Value = 1.47
$ws.Range("P15").Value = 1.5
$ws.Range("Q15").Value = 2.25
$ws.Range("R15").Value = 2.02
$ws.Range("S15").Value = 1.62
$ws.Range("T15").Value = 5.6
$ws.Range("U15").Value = 8.75
$ws.Range("V15").Value = 9.25
$ws.Range("W15").Value = 19
$ws.Range("X15").Value = 20
$ws.Range("Y15").Value = 40
$ws.Range("Z15").Value = 6.6
$ws.Range("AA15").Value = 5.9
$ws.Range("AB15").Value = 18
$ws.Range("AC15").Value = 110
$ws.Range("AE15").Value = 8.25
$ws.Range("AF15").Value = 18
$ws.Range("AG15").Value = 13
$ws.Range("AH15").Value = 55
$ws.Range("AI15").Value = 40
$ws.Range("AJ15").Value = 55
# Row 16
$ws.Range("T16").Value = 10.75
$ws.Range("X16").Value = 21
$ws.Range("AF16").Value = 13.5
